$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 629, shifting existing rows 629:660 down to 630:661
$ws.Rows.Item(629).Insert()

# Populate the new row 629 with the new data record
$ws.Cells.Item(629, 1).Value = 9
$ws.Cells.Item(629, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(629, 3).Value = "Metropolitana"
$ws.Cells.Item(629, 4).Value = 44753
$ws.Cells.Item(629, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(629, 5).Value = 13
$ws.Cells.Item(629, 6).Value = 100112040
$ws.Cells.Item(629, 7).Value = "Cilantro"
$ws.Cells.Item(629, 8).Value = "Sin especificar"
$ws.Cells.Item(629, 9).Value = "Primera"
$ws.Cells.Item(629, 10).Value = 70
$ws.Cells.Item(629, 11).Value = 15000
$ws.Cells.Item(629, 12).Value = 16000
$ws.Cells.Item(629, 13).Value = 15500
$ws.Cells.Item(629, 14).Value = "`$/docena de atados"
$ws.Cells.Item(629, 15).Value = "Región Metropolitana"
$ws.Cells.Item(629, 16).Value = 5167
$ws.Cells.Item(629, 17).Value = 3
$ws.Cells.Item(629, 18).Value = "Hortaliza"
